# Update the division-problem answers in the worksheet table.
# The document contains a single table whose populated rows are
# 1, 5, 9, 13 and 17 (each holding 5 answers, one per column).
# We address each cell directly by (row, column) to avoid any
# ambiguity caused by duplicate/overlapping old and new values.

$d = $word.ActiveDocument
$t = $d.Tables(1)

$updates = @(
    @{ Row = 1;  Col = 1; Text = "50÷8=6, 2" },
    @{ Row = 1;  Col = 2; Text = "51÷5=10, 1" },
    @{ Row = 1;  Col = 3; Text = "63÷6=10, 3" },
    @{ Row = 1;  Col = 4; Text = "98÷6=16, 2" },
    @{ Row = 1;  Col = 5; Text = "17÷6=2, 5" },

    @{ Row = 5;  Col = 1; Text = "30÷2=15, 0" },
    @{ Row = 5;  Col = 2; Text = "48÷8=6, 0" },
    @{ Row = 5;  Col = 3; Text = "17÷7=2, 3" },
    @{ Row = 5;  Col = 4; Text = "73÷3=24, 1" },
    @{ Row = 5;  Col = 5; Text = "18÷5=3, 3" },

    @{ Row = 9;  Col = 1; Text = "60÷8=7, 4" },
    @{ Row = 9;  Col = 2; Text = "26÷9=2, 8" },
    @{ Row = 9;  Col = 3; Text = "13÷9=1, 4" },
    @{ Row = 9;  Col = 4; Text = "17÷9=1, 8" },
    @{ Row = 9;  Col = 5; Text = "26÷7=3, 5" },

    @{ Row = 13; Col = 1; Text = "29÷8=3, 5" },
    @{ Row = 13; Col = 2; Text = "78÷2=39, 0" },
    @{ Row = 13; Col = 3; Text = "56÷5=11, 1" },
    @{ Row = 13; Col = 4; Text = "62÷7=8, 6" },
    @{ Row = 13; Col = 5; Text = "34÷7=4, 6" },

    @{ Row = 17; Col = 1; Text = "10÷6=1, 4" },
    @{ Row = 17; Col = 2; Text = "14÷9=1, 5" },
    @{ Row = 17; Col = 3; Text = "87÷8=10, 7" },
    @{ Row = 17; Col = 4; Text = "45÷9=5, 0" },
    @{ Row = 17; Col = 5; Text = "81÷6=13, 3" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
